# Updates cryptos list values (price and 1h volume change) per latest scrape.
# Row 17/18 also swap Chainlink/WrappedBTC positions (rank change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "26.985.31"
$ws.Cells.Item(2, 5).Value = "  -0.08%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.677.70"
$ws.Cells.Item(3, 5).Value = "  +0.45%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.13%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'215.09"
$ws.Cells.Item(5, 5).Value = "  -0.42%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.518"
$ws.Cells.Item(6, 5).Value = "  +1.16%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.19%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.0619"
$ws.Cells.Item(9, 5).Value = "  +0.33%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'20.33"
$ws.Cells.Item(10, 5).Value = "  +1.20%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0886"
$ws.Cells.Item(11, 5).Value = "  -0.55%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.915.42"
$ws.Cells.Item(12, 5).Value = "  +0.46%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.681.41"
$ws.Cells.Item(13, 5).Value = "  +0.63%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.26%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.528"
$ws.Cells.Item(15, 5).Value = "  +1.56%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'65.70"
$ws.Cells.Item(16, 5).Value = "  -0.06%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "27.006.87"
$ws.Cells.Item(17, 5).Value = "  -0.04%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).Value = "'8.19"
$ws.Cells.Item(18, 5).Value = "  +6.11%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'235.59"
$ws.Cells.Item(19, 5).Value = "  +0.13%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.17%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.03%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'4.44"

# Row 23
$ws.Cells.Item(23, 4).Value = "'9.18"
$ws.Cells.Item(23, 5).Value = "  -0.55%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'2.16"
$ws.Cells.Item(24, 5).Value = "  -2.89%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'145.95"
$ws.Cells.Item(25, 5).Value = "  +0.35%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'7.22"
$ws.Cells.Item(26, 5).Value = "  +0.78%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +1.35%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'0.113"
$ws.Cells.Item(28, 5).Value = "  -1.43%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.16%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'0.0497"
$ws.Cells.Item(30, 5).Value = "  -0.16%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.65%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.14%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.478.98"
$ws.Cells.Item(33, 5).Value = "  +1.85%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +1.34%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +4.88%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.22%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.582"
$ws.Cells.Item(37, 5).Value = "  +2.33%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.37%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.903"
$ws.Cells.Item(39, 5).Value = "  +1.32%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'5.83"
$ws.Cells.Item(40, 5).Value = "  -3.91%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.51%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.12%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +1.41%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'67.36"
$ws.Cells.Item(44, 5).Value = "  +2.50%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.819.10"
$ws.Cells.Item(45, 5).Value = "  +0.12%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.782"
$ws.Cells.Item(46, 5).Value = "  +0.27%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'90.39"
$ws.Cells.Item(47, 5).Value = "  +0.13%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.81%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.57%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +1.55%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +0.18%  "
